# Auto update Excel log
# Appends new motion-sensor log rows (97-103) to the "mmWave" worksheet,
# mirroring the same Date / Timestamp / Hour / Location / Value / Status
# layout used by the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-02-01", "18:50:53", "18:00", "Living Room", "NO_MOTION_DETECTED", "Inactive"),
    @("2026-02-01", "18:51:00", "18:00", "Living Room", "PRESENCE_DETECTED",  "Active"),
    @("2026-02-01", "18:51:10", "18:00", "Living Room", "PRESENCE_DETECTED",  "Active"),
    @("2026-02-01", "18:51:21", "18:00", "Living Room", "PRESENCE_DETECTED",  "Active"),
    @("2026-02-01", "18:51:31", "18:00", "Living Room", "PRESENCE_DETECTED",  "Active"),
    @("2026-02-01", "18:51:42", "18:00", "Living Room", "PRESENCE_DETECTED",  "Active"),
    @("2026-02-01", "18:51:52", "18:00", "Living Room", "PRESENCE_DETECTED",  "Active")
)

$startRow = 97
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    # Column A holds a plain "Date" string (e.g. "2026-02-01"). A bare
    # .Value assignment gets auto-recognized as a real date and converted
    # to a serial number, so force Text formatting for the write, then
    # drop the formatting again so the cell ends up identical to its
    # plain-text siblings (no leftover style reference).
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]
    $cellA.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
    $ws.Cells.Item($r, 6).Value = $values[5]
}
